# Aligning 2D NMR Spectra Part 2 - apply commit "Add files via upload"
#
# 1) In the mlrMBO paragraph, delete the sentence "There is a nice
#    introductory vignette." (this removes the hyperlinked text
#    "introductory vignette" along with the surrounding plain text).
# 2) In the "lhs package" sentence, unlink (but keep) the hyperlinked
#    text "lhs package" - remove the hyperlink, keep the run/text.
# 3) In the "mlrMBO provides many options for the surrogate model"
#    sentence, remove the word "many " and unlink+remove the
#    hyperlinked text "options" so it reads
#    "mlrMBO provides for the surrogate model."

$d = $word.ActiveDocument

# --- Edit 1: remove " There is a nice introductory vignette." ---
$found1 = $d.Content.Find.Execute(
    "There is a nice introductory vignette.",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 2)

# --- Edit 2: unlink "lhs package" (keep text, drop hyperlink) ---
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "lhs package") {
        $h.Delete()
        break
    }
}

# --- Edit 3: "provides many <options>" -> "provides " (then drop the
#     now-unused "options" hyperlink's remaining text too) ---
$found3 = $d.Content.Find.Execute(
    "provides many options for",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "provides for",
    2)

Write-Output "find1=$found1 find3=$found3"
